$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new LF07 / Krank / Recherche entries ---
$ws.Range("L2").Value = "LF07 Stunde"
$ws.Range("M2").Value = "Krank"
$ws.Range("N2").Value = "Recherche wie eine readme Datei erstellt wird"

# --- Row 3: new LF07 / Bugs / Powerpoint entries ---
$ws.Range("L3").Value = "LF07 Stunde"
$ws.Range("M3").Value = "Bugs gefixt"
$ws.Range("N3").Value = "Powerpoint angefangen"

# --- Row 4: newly filled-in "Krank" cells plus LF07 / Java Doc entries ---
$ws.Range("C4").Value = "Krank"
$ws.Range("F4").Value = "Versuch von Bug fixes"
$ws.Range("G4").Value = "Krank"
$ws.Range("H4").Value = "Krank"
$ws.Range("I4").Value = "Krank"
$ws.Range("L4").Value = "LF07 stunde"
$ws.Range("M4").Value = "Beschreibungen von Gegenständen in Java Doc erstellt"
$ws.Range("N4").Value = "Finalisierung Java Doc"

# --- Column widths (F widened, L/M/N newly sized for the added columns) ---
$ws.Columns.Item(6).ColumnWidth = 24.42578125
$ws.Columns.Item(12).ColumnWidth = 16
$ws.Columns.Item(13).ColumnWidth = 49.7109375
$ws.Columns.Item(14).ColumnWidth = 42.140625

# --- View state: scroll window right to column K, select I4 ---
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("I4").Select()
